$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.704.03'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '1.605.15'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  +0.52%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.80'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.09'
$ws.Range('E8').Value = '  +4.72%  '
$ws.Range('E9').Value = '  +0.91%  '
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0910'
$ws.Range('E11').Value = '  -0.51%  '
$ws.Range('D12').Value = '1.834.93'
$ws.Range('D13').Value = '1.612.51'
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('E14').Value = '  +3.75%  '
$ws.Range('D15').Value = '29.710.59'
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.08'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '241.34'
$ws.Range('E18').Value = '  -2.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.87'
$ws.Range('E19').Value = '  +3.38%  '
$ws.Range('D20').Value = '0.0₃0697'
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('E22').Value = '  -0.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.40'
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.51'
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.46'
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.46'
$ws.Range('E28').Value = '  +0.79%  '
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('E30').Value = '  +1.29%  '
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('E32').Value = '  -0.17%  '
$ws.Range('E33').Value = '  +2.09%  '
$ws.Range('D34').Value = '1.425.76'
$ws.Range('E34').Value = '  -0.92%  '
$ws.Range('E35').Value = '  +3.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.90'
$ws.Range('E36').Value = '  +1.81%  '
$ws.Range('E37').Value = '  -1.93%  '
$ws.Range('E38').Value = '  -0.24%  '
$ws.Range('E39').Value = '  +1.57%  '
$ws.Range('E40').Value = '  +2.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '56.81'
$ws.Range('E41').Value = '  +3.09%  '
$ws.Range('E42').Value = '  +6.12%  '
$ws.Range('E43').Value = '  +2.06%  '
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('E45').Value = '  +0.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '66.31'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.982'
$ws.Range('E47').Value = '  +17.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.39'
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').Value = '1.743.44'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.62'
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range('E51').Value = '  +4.05%  '
